$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.128.09"
$ws.Range("E2").Value = "  -1.19%  "
$ws.Range("D3").Value = "2.427.68"
$ws.Range("E3").Value = "  -1.77%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.00"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "89.40"
$ws.Range("E6").Value = "  -3.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.539"
$ws.Range("E7").Value = "  -2.37%  "
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.498"
$ws.Range("E9").Value = "  -3.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0836"
$ws.Range("E10").Value = "  -1.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "32.05"
$ws.Range("E11").Value = "  -2.93%  "
$ws.Range("E12").Value = "  -1.53%  "
$ws.Range("D13").Value = "2.803.70"
$ws.Range("E13").Value = "  -1.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.73"
$ws.Range("E14").Value = "  -2.41%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.64"
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("D16").Value = "2.460.88"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.777"
$ws.Range("E17").Value = "  -1.90%  "
$ws.Range("D18").Value = "41.077.82"
$ws.Range("E18").Value = "  -1.22%  "
$ws.Range("D19").Value = "0.0₃0926"
$ws.Range("E19").Value = "  -2.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.29"
$ws.Range("E20").Value = "  -2.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.31"
$ws.Range("E21").Value = "  +1.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.05"
$ws.Range("E22").Value = "  -2.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.20"
$ws.Range("E23").Value = "  -1.79%  "
$ws.Range("E24").Value = "  -1.66%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("E26").Value = "  -2.81%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.14"
$ws.Range("E27").Value = "  -2.09%  "
$ws.Range("E28").Value = "  -2.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.63"
$ws.Range("E29").Value = "  -2.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.55"
$ws.Range("E30").Value = "  -4.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "158.42"
$ws.Range("E31").Value = "  -1.69%  "
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.26"
$ws.Range("E33").Value = "  -4.63%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0746"
$ws.Range("E34").Value = "  -2.61%  "
$ws.Range("E35").Value = "  -3.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.04"
$ws.Range("E36").Value = "  -1.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.94"
$ws.Range("E37").Value = "  +1.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.114"
$ws.Range("E38").Value = "  -1.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.79"
$ws.Range("E39").Value = "  -3.47%  "
$ws.Range("E40").Value = "  -2.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.89"
$ws.Range("E41").Value = "  -2.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.33"
$ws.Range("E42").Value = "  -5.17%  "
$ws.Range("D43").Value = "1.995.36"
$ws.Range("E43").Value = "  +0.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.53"
$ws.Range("E44").Value = "  -2.30%  "
$ws.Range("E45").Value = "  -3.05%  "
$ws.Range("E46").Value = "  -2.92%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.54"
$ws.Range("E47").Value = "  +3.73%  "
$ws.Range("D48").Value = "2.666.39"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "94.75"
$ws.Range("E49").Value = "  -2.76%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.40"
$ws.Range("E50").Value = "  -1.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.91"
$ws.Range("E51").Value = "  -0.66%  "
